$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simplify combined taxonomy/method values that included "&2M" suffix down to
# their base value, and turn the lone "2M" entry into "NONE" (adding methods
# to get bacterium taxonomy by its ID).
$ws.Range("C8").Value = "CL"
$ws.Range("B9").Value = "NONE"
$ws.Range("C9").Value = "SCL"
$ws.Range("E9").Value = "OL"
$ws.Range("F9").Value = "OL"
$ws.Range("G9").Value = "CL"
$ws.Range("D11").Value = "CL"
$ws.Range("F11").Value = "SCL"
$ws.Range("G11").Value = "CL"
$ws.Range("D27").Value = "CL"
$ws.Range("D29").Value = "CL"
$ws.Range("G29").Value = "CL"
$ws.Range("G31").Value = "CL"
$ws.Range("H31").Value = "CL"
$ws.Range("G33").Value = "CL"
$ws.Range("E40").Value = "SCL"
$ws.Range("D55").Value = "CL"
$ws.Range("E55").Value = "SCL"
$ws.Range("G55").Value = "CL"
$ws.Range("B61").Value = "SCL"
$ws.Range("C61").Value = "SCL"
$ws.Range("E62").Value = "SCL"

# Update the active selection to match the final saved view state.
$ws.Range("B10").Select()
